# Hungary NB I workbook update (09-04-2024)
# A new match (Man Utd vs Man City, played 2024-04-09) is inserted as a new
# data row right before the existing row that used to be row 162 (id 160,
# match 7939469). All rows from the old row 162 onward shift down by one,
# and a handful of odds cells on the shifted rows were re-scraped with
# slightly different values, matching the source diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new blank row at sheet row 162 -------------------------------
# This pushes the former rows 162-166 down to 163-167, carrying their
# existing values/formats with them automatically.
$ws.Rows(162).Insert()

# The freshly inserted row has no number formatting yet; clone it from the
# row immediately below (which now holds what used to be row 162), so the
# date cell (E162) keeps the date/time display format etc.
$ws.Range("A163:AC163").Copy()
$ws.Range("A162:AC162").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Populate the new row 162 with the new match -----------------------
$ws.Cells.Item(162, 1).Value = 160
$ws.Cells.Item(162, 2).Value = 8074943
$ws.Cells.Item(162, 3).Value = "Hungary NB I"
$ws.Cells.Item(162, 4).Value = "Hungary NB I"
$ws.Cells.Item(162, 5).Value = 45391.53819444445
$ws.Cells.Item(162, 6).Value = "Man Utd"
$ws.Cells.Item(162, 7).Value = "Man City"
$ws.Cells.Item(162, 8).Value = 2
$ws.Cells.Item(162, 9).Value = 2
$ws.Cells.Item(162, 10).Value = "D"
$ws.Cells.Item(162, 11).Value = 8.5
$ws.Cells.Item(162, 12).Value = 3.8
$ws.Cells.Item(162, 13).Value = 1.5
$ws.Cells.Item(162, 14).Value = 8.5
$ws.Cells.Item(162, 15).Value = 3.8
$ws.Cells.Item(162, 16).Value = 1.5
$ws.Cells.Item(162, 17).Value = 1
$ws.Cells.Item(162, 18).Value = 1.95
$ws.Cells.Item(162, 19).Value = 1.9
$ws.Cells.Item(162, 20).Value = 2.5
$ws.Cells.Item(162, 21).Value = 1.85
$ws.Cells.Item(162, 22).Value = 2
$ws.Cells.Item(162, 23).Value = -1
$ws.Cells.Item(162, 24).Value = 2.8
$ws.Cells.Item(162, 25).Value = -1
$ws.Cells.Item(162, 26).Value = 0.95
$ws.Cells.Item(162, 27).Value = -1
$ws.Cells.Item(162, 28).Value = 0.8500000000000001
$ws.Cells.Item(162, 29).Value = -1

# --- Re-scraped odds on the rows that shifted down ----------------------
# Column A is a plain row sequence number (0-based, = sheet row - 2), not
# part of the shifted match data, so restore it on every row the insert
# touched (the insert operation dragged the old values down along with
# everything else, which is wrong for this particular column).
$ws.Cells.Item(163, 1).Value = 161
$ws.Cells.Item(164, 1).Value = 162
$ws.Cells.Item(165, 1).Value = 163
$ws.Cells.Item(166, 1).Value = 164
$ws.Cells.Item(167, 1).Value = 165

# Row 163 (was row 162, id 160, match 7939469)
$ws.Cells.Item(163, 14).Value = 1.285
$ws.Cells.Item(163, 16).Value = 9
$ws.Cells.Item(163, 17).Value = -1.75
$ws.Cells.Item(163, 18).Value = 2.025
$ws.Cells.Item(163, 19).Value = 1.825
$ws.Cells.Item(163, 20).Value = 3
$ws.Cells.Item(163, 21).Value = 1.85
$ws.Cells.Item(163, 22).Value = 2

# Row 164 (was row 163, id 161, match 6818359)
$ws.Cells.Item(164, 14).Value = 2.5
$ws.Cells.Item(164, 16).Value = 2.7
$ws.Cells.Item(164, 18).Value = 1.8
$ws.Cells.Item(164, 19).Value = 2.05

# Row 165 (was row 164, id 162, match 6818358)
$ws.Cells.Item(165, 14).Value = 1.571
$ws.Cells.Item(165, 15).Value = 3.75
$ws.Cells.Item(165, 16).Value = 6.5
$ws.Cells.Item(165, 17).Value = -1
$ws.Cells.Item(165, 18).Value = 2
$ws.Cells.Item(165, 19).Value = 1.85
$ws.Cells.Item(165, 21).Value = 1.925
$ws.Cells.Item(165, 22).Value = 1.925

# Row 166 (was row 165, id 163, match 6818360)
$ws.Cells.Item(166, 14).Value = 1.65
$ws.Cells.Item(166, 15).Value = 3.6
$ws.Cells.Item(166, 16).Value = 5.5
$ws.Cells.Item(166, 17).Value = -0.75
$ws.Cells.Item(166, 18).Value = 1.8
$ws.Cells.Item(166, 19).Value = 2.05

# Row 167 (was row 166, id 164, match 6818361)
$ws.Cells.Item(167, 18).Value = 1.95
$ws.Cells.Item(167, 19).Value = 1.9
